# 10Th - MB for single stock and added new group
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column B (the most-recent "week" column),
# shifting the existing week columns (old B, old C) one slot to the right.
$ws.Columns.Item(2).Insert()

# Match the width of the neighboring "week" column.
$ws.Columns.Item(2).ColumnWidth = 30.83203125

# New header for the freshly inserted "this week" column.
$ws.Range("B1").Value = "Jun_27"

# Every analyst row gets "UN" (unchanged) in the new column B, same as
# the rest of the sheet's default/unchanged marker.
$lastRow = 27
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
}

# Two new research firms tracked starting this week.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"

$ws.Range("F8").Select()
